# Add a new slide (Title and Content layout = index 2 / ppLayoutText) after
# the existing Timeline slide, documenting the source the timeline was
# built from ("added source link to how the timeline was created").

$p = $ppt.ActivePresentation

# ppLayoutText (2) -> "Title and Content" slide layout, inserted as slide 2.
$newSlide = $p.Slides.Add(2, 2)

# Leave the Title placeholder empty (matches the authored slide) and put the
# hyperlinked source text into the body/content placeholder.
$contentShape = $newSlide.Shapes.Item(2)
$contentShape.TextFrame.TextRange.Text = "Presentationgo.com"

$linkRange = $contentShape.TextFrame.TextRange
$hyperlink = $linkRange.ActionSettings.Item(1).Hyperlink
$hyperlink.Address = "http://presentationgo.com/"
$hyperlink.ScreenTip = "http://presentationgo.com/"
